# Adds two new columns, I ("I0") and J ("IF"), to the stats sheet,
# mirroring the existing header style and filling in the per-row data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers "I0" and "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered alignment) used by the
# other header cells (e.g. H1) onto the two new header cells, without
# touching the values we just set.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (2-39): values for columns I and J ---
$values = @(
    ,(1,4)
    ,(9,9)
    ,(9,9)
    ,(9,9)
    ,(9,9)
    ,(7,8)
    ,(6,7)
    ,(7,7)
    ,(6,7)
    ,(9,9)
    ,(7,8)
    ,(7,7)
    ,(4,6)
    ,(3,4)
    ,(7,7)
    ,(8,9)
    ,(9,9)
    ,(7,7)
    ,(6,7)
    ,(7,8)
    ,(4,4)
    ,(7,7)
    ,(9,9)
    ,(6,7)
    ,(6,7)
    ,(9,9)
    ,(6,6)
    ,(8,9)
    ,(6,6)
    ,(8,8)
    ,(7,7)
    ,(5,5)
    ,(9,9)
    ,(6,6)
    ,(6,6)
    ,(5,5)
    ,(5,5)
    ,(5,5)
)

$rowCount = $values.Count
$arr = New-Object 'object[,]' $rowCount,2
for ($r = 0; $r -lt $rowCount; $r++) {
    $arr[$r,0] = $values[$r][0]
    $arr[$r,1] = $values[$r][1]
}

$ws.Range("I2:J39").Value = $arr

Write-Host "Added I0/IF columns for $rowCount data rows."
